# Update gh-pages output data at 456a3b4
# Sheet "展览" (Exhibition) - column F (想去人数 / "want to go" count)
$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value  = 1778
$wsExhibit.Range("F5").Value  = 3326
$wsExhibit.Range("F6").Value  = 1052
$wsExhibit.Range("F7").Value  = 2206
$wsExhibit.Range("F8").Value  = 2119
$wsExhibit.Range("F9").Value  = 1108
$wsExhibit.Range("F10").Value = 608
$wsExhibit.Range("F13").Value = 396
$wsExhibit.Range("F15").Value = 43
$wsExhibit.Range("F16").Value = 95
$wsExhibit.Range("F17").Value = 213
$wsExhibit.Range("F18").Value = 1586
$wsExhibit.Range("F19").Value = 633
$wsExhibit.Range("F20").Value = 725
$wsExhibit.Range("F21").Value = 605
$wsExhibit.Range("F22").Value = 12251
$wsExhibit.Range("F23").Value = 12295
$wsExhibit.Range("F24").Value = 910
$wsExhibit.Range("F27").Value = 36
$wsExhibit.Range("F28").Value = 21
$wsExhibit.Range("F29").Value = 368
$wsExhibit.Range("F32").Value = 200
$wsExhibit.Range("F33").Value = 586

# Sheet "演出" (Performance) - column F
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F7").Value = 28

# Sheet "全部类型" (All types) - column F
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value  = 1778
$wsAll.Range("F6").Value  = 3326
$wsAll.Range("F7").Value  = 1052
$wsAll.Range("F8").Value  = 2206
$wsAll.Range("F9").Value  = 2119
$wsAll.Range("F10").Value = 1108
$wsAll.Range("F11").Value = 608
$wsAll.Range("F14").Value = 396
$wsAll.Range("F17").Value = 43
$wsAll.Range("F19").Value = 95
$wsAll.Range("F21").Value = 213
$wsAll.Range("F22").Value = 1586
$wsAll.Range("F23").Value = 633
$wsAll.Range("F24").Value = 725
$wsAll.Range("F25").Value = 605
$wsAll.Range("F26").Value = 12251
$wsAll.Range("F27").Value = 12295
$wsAll.Range("F28").Value = 910
$wsAll.Range("F31").Value = 36
$wsAll.Range("F32").Value = 21
$wsAll.Range("F33").Value = 368
$wsAll.Range("F38").Value = 200
$wsAll.Range("F39").Value = 586
$wsAll.Range("F40").Value = 28
